$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block of weekly "Ajo" price records (rows 80-104) is shifted down
# by two rows (to rows 82-106) to make room for a new week's worth of data,
# which is inserted at the (now vacated) top of the block: rows 80-81.
$ws.Range("A80:R81").EntireRow.Insert()

# New row 80: Ajo / Chino / Primera, $/caja 10 kilos, origin China
$ws.Range("A80").Value = 9
$ws.Range("B80").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C80").Value = "Metropolitana"
$ws.Range("D80").Value = 44452
$ws.Range("E80").Value = 13
$ws.Range("F80").Value = 100112003
$ws.Range("G80").Value = "Ajo"
$ws.Range("H80").Value = "Chino"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 430
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 15500
$ws.Range("M80").Value = 15250
$ws.Range("N80").Value = "$/caja 10 kilos"
$ws.Range("O80").Value = "China"
$ws.Range("P80").Value = 1525
$ws.Range("Q80").Value = 10
$ws.Range("R80").Value = "Hortaliza"

# New row 81: Ajo / Chino / Primera, $/malla 10 kilos, origin China
$ws.Range("A81").Value = 9
$ws.Range("B81").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value = 44452
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = 100112003
$ws.Range("G81").Value = "Ajo"
$ws.Range("H81").Value = "Chino"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 340
$ws.Range("K81").Value = 15500
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = 15750
$ws.Range("N81").Value = "$/malla 10 kilos"
$ws.Range("O81").Value = "China"
$ws.Range("P81").Value = 1575
$ws.Range("Q81").Value = 10
$ws.Range("R81").Value = "Hortaliza"
